# Insert a new data row at row 17 ("LO / CAP_BND / 2045 / 16 / ELE_NEW_WIND-ON")
# in the "Nowe moce na morzu" (offshore wind) block, pushing every row from the
# old row 17 onward down by one (old row 17 -> 18, ... old row 45 -> 46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row above the old row 17, shifting rows 17:45 down to 18:46.
$ws.Rows(17).Insert()

# Match the row height used by the data rows in this block (row 16 is a sibling
# data row with the same formatting that should carry through to the new row).
$ws.Rows(17).RowHeight = $ws.Rows(16).RowHeight

# Fill in the new row's data.
$ws.Range("B17").Value = "LO"
$ws.Range("C17").Value = "CAP_BND"
$ws.Range("D17").Value = 2045
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = "ELE_NEW_WIND-ON"

# Match the updated selection recorded in the saved workbook.
$ws.Range("J18").Select()
